$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2026-02-28 -> 2026-03-01, i.e. 46081 -> 46082) for every data row.
$ws.Range("C2:C247").Value = 46082
